$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header (day numbers -> day names in Portuguese)
$ws.Range("B1").Value = "segunda"
$ws.Range("C1").Value = "terça"
$ws.Range("D1").Value = "quarta"
$ws.Range("E1").Value = "quinta"
$ws.Range("F1").Value = "sexta"

# Column A - period index -> time of day
$ws.Range("A2").Value = "7:00"
$ws.Range("A3").Value = "7:50"
$ws.Range("A4").Value = "8:40"
$ws.Range("A5").Value = "9:30"
$ws.Range("A6").Value = "10:40"
$ws.Range("A7").Value = "11:30"
$ws.Range("A8").Value = "13:00"
$ws.Range("A9").Value = "13:50"
$ws.Range("A10").Value = "14:40"
$ws.Range("A11").Value = "15:30"
$ws.Range("A12").Value = "16:40"
$ws.Range("A13").Value = "17:30"

# Content shifts within the timetable grid
$ws.Range("E4").Value = "-"
$ws.Range("E5").Value = "Circuitos Elétricos 2"
$ws.Range("B6").Value = "Desenho Técnico"
$ws.Range("F9").Value = "-"
$ws.Range("B10").Value = "EAP"
$ws.Range("D10").Value = "EAP"
$ws.Range("E11").Value = "-"
